$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author's edit also cleared the stale cell selection (A14:C15) that was
# saved in the sheet view, leaving the default single-cell selection at A1.
$ws.Range("A1").Select()

# Add a new "2021" data column (column I) mirroring the existing
# 2016-2020 columns (D-H) for each data row (4-25).
$ws.Range("I4").Value = 2021
$ws.Range("I5").Value = 48.5
$ws.Range("I7").Value = 48.8
$ws.Range("I8").Value = 48.2
$ws.Range("I10").Value = 58.2
$ws.Range("I11").Value = 42.4
$ws.Range("I12").Value = 40.7
$ws.Range("I14").Value = 41.5
$ws.Range("I15").Value = 52.6
$ws.Range("I17").Value = 67.1
$ws.Range("I18").Value = 62
$ws.Range("I19").Value = 46.9
$ws.Range("I20").Value = 55.8
$ws.Range("I21").Value = 42.7
$ws.Range("I22").Value = 48.3
$ws.Range("I23").Value = 39.7
$ws.Range("I24").Value = 38.1
$ws.Range("I25").Value = 44.7

# Rows 6, 9, 13 and 16 are section headers with no numeric value in columns
# D-H, just the shared number formatting - replicate that for column I too.

# Copy the formatting (number format/borders/etc.) from column H, row by row,
# onto the new column I so every new cell matches its row's existing style.
$ws.Range("H4:H25").Copy()
$ws.Range("I4:I25").PasteSpecial(-4122)
